$wb = $excel.ActiveWorkbook

# Sheet ALC, row 126
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 41740
$ws.Range("J126").Value = 41740
$ws.Range("L126").Value = 41740
$ws.Range("N126").Value = -51620

# Sheet ALC, row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 720.73334
$ws.Range("I129").Value = 525.9167
$ws.Range("J129").Value = 1500
$ws.Range("K129").Value = 1577.7501
$ws.Range("L129").Value = 4500
$ws.Range("M129").Value = 3422.2499
$ws.Range("N129").Value = -14500

# Sheet ALC, row 130
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 38838.75
$ws.Range("J130").Value = 38838.75
$ws.Range("L130").Value = 38838.75
$ws.Range("N130").Value = -48878.75

# Sheet ALC, row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 60767
$ws.Range("J133").Value = 60767
$ws.Range("L133").Value = 60767
$ws.Range("N133").Value = -70887

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 58825170
$ws.Range("I102").Value = 1625.4546
$ws.Range("J102").Value = 166668340
$ws.Range("K102").Value = 1625.4546
$ws.Range("L102").Value = 166668340
$ws.Range("M102").Value = -3.454600000000028
$ws.Range("N102").Value = -166671584

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 13189.895
$ws.Range("I122").Value = 18548.334
$ws.Range("J122").Value = 4004
$ws.Range("K122").Value = 55645.00199999999
$ws.Range("L122").Value = 12012
$ws.Range("M122").Value = -53195.00199999999
$ws.Range("N122").Value = -16912

# Sheet ARM, row 125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 180033840
$ws.Range("J125").Value = 180033840
$ws.Range("L125").Value = 180033840
$ws.Range("N125").Value = -180043680

# Sheet ARM, row 129
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999

# Sheet ARM, row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 28240.2
$ws.Range("J133").Value = 28240.2
$ws.Range("L133").Value = 28240.2
$ws.Range("N133").Value = -33300.2

# Sheet ARM, row 134
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 92333
$ws.Range("J134").Value = 92333
$ws.Range("L134").Value = 92333
$ws.Range("N134").Value = -102473

# Sheet BSM, row 26
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 2406.8333
$ws.Range("I26").Value = 2406.8333
$ws.Range("K26").Value = 2406.8333
$ws.Range("M26").Value = -2114.8333

# Sheet BSM, row 50
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 20500
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 20500
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 20500
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -21648

# Sheet BSM, row 51
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 14545
$ws.Range("J51").Value = 14545
$ws.Range("L51").Value = 14545
$ws.Range("N51").Value = -15527

# Sheet BSM, row 96
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 14618.667
$ws.Range("I96").Value = 14618.667
$ws.Range("K96").Value = 14618.667
$ws.Range("M96").Value = -11872.667

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1674.5454
$ws.Range("I99").Value = 1058.8889
$ws.Range("J99").Value = 4445
$ws.Range("K99").Value = 1058.8889
$ws.Range("L99").Value = 4445
$ws.Range("M99").Value = 439.1111000000001
$ws.Range("N99").Value = -7441

# Sheet BSM, row 125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Sheet BSM, row 126
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 27981
$ws.Range("J126").Value = 27981
$ws.Range("L126").Value = 27981
$ws.Range("N126").Value = -37861

# Sheet BSM, row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 42922.855
$ws.Range("J132").Value = 42922.855
$ws.Range("L132").Value = 42922.855
$ws.Range("N132").Value = -53042.855

# Sheet CUL, row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 148.33333
$ws.Range("I26").Value = 125
$ws.Range("J26").Value = 195
$ws.Range("K26").Value = 375
$ws.Range("L26").Value = 585
$ws.Range("M26").Value = -87
$ws.Range("N26").Value = -1161

# Sheet CUL, row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 252691.19
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 252691.19
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 758073.5700000001
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -760693.5700000001

# Sheet GSM, row 46
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 9223.777
$ws.Range("J46").Value = 10419
$ws.Range("L46").Value = 10419
$ws.Range("N46").Value = -10731

# Sheet GSM, row 62
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Sheet GSM, row 65
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# Sheet GSM, row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1400
$ws.Range("I97").Value = 1333.3334
$ws.Range("J97").Value = 1466.6666
$ws.Range("K97").Value = 1333.3334
$ws.Range("L97").Value = 1466.6666
$ws.Range("M97").Value = -837.3334
$ws.Range("N97").Value = -2458.6666

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2726.7273
$ws.Range("I122").Value = 2570.7144
$ws.Range("J122").Value = 2999.75
$ws.Range("K122").Value = 7712.1432
$ws.Range("L122").Value = 8999.25
$ws.Range("M122").Value = -5262.1432
$ws.Range("N122").Value = -13899.25

# Sheet GSM, row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 16374.421
$ws.Range("I126").Value = 3050
$ws.Range("J126").Value = 19927.6
$ws.Range("K126").Value = 9150
$ws.Range("L126").Value = 59782.8
$ws.Range("M126").Value = -6680
$ws.Range("N126").Value = -64722.8

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2566.5454
$ws.Range("I132").Value = 2507.36
$ws.Range("J132").Value = 2751.5
$ws.Range("K132").Value = 7522.08
$ws.Range("L132").Value = 8254.5
$ws.Range("M132").Value = -4992.08
$ws.Range("N132").Value = -13314.5

# Sheet GSM, row 135
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 46950.527
$ws.Range("J135").Value = 46950.527
$ws.Range("L135").Value = 46950.527
$ws.Range("N135").Value = -57090.527

# Sheet LTW, row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2114.2856
$ws.Range("I100").Value = 1633.4445
$ws.Range("J100").Value = 2979.8
$ws.Range("K100").Value = 1633.4445
$ws.Range("L100").Value = 2979.8
$ws.Range("M100").Value = -1092.4445
$ws.Range("N100").Value = -4061.8

# Sheet LTW, row 108
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 28643.4
$ws.Range("J108").Value = 28643.4
$ws.Range("L108").Value = 28643.4
$ws.Range("N108").Value = -36323.4

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5380
$ws.Range("I122").Value = 4571.4287
$ws.Range("J122").Value = 7266.6665
$ws.Range("K122").Value = 13714.2861
$ws.Range("L122").Value = 21799.9995
$ws.Range("M122").Value = -11264.2861
$ws.Range("N122").Value = -26699.9995

# Sheet LTW, row 125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 41994
$ws.Range("J125").Value = 41994
$ws.Range("L125").Value = 41994
$ws.Range("N125").Value = -51834

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3642.0476
$ws.Range("I132").Value = 3644
$ws.Range("J132").Value = 3639.9
$ws.Range("K132").Value = 10932
$ws.Range("L132").Value = 10919.7
$ws.Range("M132").Value = -8402
$ws.Range("N132").Value = -15979.7

# Sheet LTW, row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 45617.332
$ws.Range("J133").Value = 45617.332
$ws.Range("L133").Value = 45617.332
$ws.Range("N133").Value = -50677.332

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4712.5
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 4783.3335
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 4783.3335
$ws.Range("M62").Value = -3876
$ws.Range("N62").Value = -6031.3335

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4712.5
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 4783.3335
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 23916.6675
$ws.Range("M65").Value = -19380
$ws.Range("N65").Value = -30156.6675

# Sheet WVR, row 108
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 32500
$ws.Range("J108").Value = 32500
$ws.Range("L108").Value = 32500
$ws.Range("N108").Value = -40180

# Sheet WVR, row 128
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 230761
$ws.Range("J128").Value = 230761
$ws.Range("L128").Value = 230761
$ws.Range("N128").Value = -240721

# Sheet WVR, row 130
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 43729.332
$ws.Range("J130").Value = 43729.332
$ws.Range("L130").Value = 43729.332
$ws.Range("N130").Value = -53769.332

# Sheet WVR, row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 70522
$ws.Range("J135").Value = 70522
$ws.Range("L135").Value = 70522
$ws.Range("N135").Value = -80662
